$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 612
$ws.Cells.Item(3, 6).Value = 10609
$ws.Cells.Item(5, 6).Value = 97
$ws.Cells.Item(6, 6).Value = 663
$ws.Cells.Item(7, 6).Value = 141
$ws.Cells.Item(8, 6).Value = 12640
$ws.Cells.Item(9, 6).Value = 13050
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(14, 6).Value = 110
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(17, 6).Value = 1441
$ws.Cells.Item(19, 6).Value = 2029
$ws.Cells.Item(20, 6).Value = 1048
$ws.Cells.Item(21, 6).Value = 1582
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(25, 6).Value = 738
$ws.Cells.Item(26, 6).Value = 3042
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(28, 6).Value = 2079
$ws.Cells.Item(29, 6).Value = 8
$ws.Cells.Item(30, 6).Value = 109
$ws.Cells.Item(31, 6).Value = 1695
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(34, 6).Value = 53
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(36, 6).Value = 3792
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(39, 6).Value = 133
$ws.Cells.Item(42, 6).Value = 2127
$ws.Cells.Item(43, 6).Value = 36
$ws.Cells.Item(47, 6).Value = 36
$ws.Cells.Item(48, 6).Value = 0

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(12, 6).Value = 5
$ws.Cells.Item(15, 6).Value = 7
$ws.Cells.Item(19, 6).Value = 14
$ws.Cells.Item(20, 6).Value = 6
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(28, 6).Value = 64
$ws.Cells.Item(29, 6).Value = 8
$ws.Cells.Item(30, 6).Value = 2

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 6574

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 612
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(5, 6).Value = 663
$ws.Cells.Item(7, 6).Value = 84
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(9, 6).Value = 13050
$ws.Cells.Item(10, 6).Value = 40
$ws.Cells.Item(11, 6).Value = 1320
$ws.Cells.Item(12, 6).Value = 1299
$ws.Cells.Item(13, 6).Value = 5519
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(15, 6).Value = 366
$ws.Cells.Item(16, 6).Value = 195
$ws.Cells.Item(17, 6).Value = 1441
$ws.Cells.Item(18, 6).Value = 365
$ws.Cells.Item(19, 6).Value = 2029
$ws.Cells.Item(20, 6).Value = 1048
$ws.Cells.Item(21, 6).Value = 1582
$ws.Cells.Item(22, 6).Value = 885
$ws.Cells.Item(23, 6).Value = 514
$ws.Cells.Item(24, 6).Value = 3042
$ws.Cells.Item(25, 6).Value = 5
$ws.Cells.Item(27, 6).Value = 2079
$ws.Cells.Item(29, 6).Value = 7
$ws.Cells.Item(30, 6).Value = 1695
$ws.Cells.Item(32, 6).Value = 1008
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(35, 6).Value = 106
$ws.Cells.Item(37, 6).Value = 3792
$ws.Cells.Item(38, 6).Value = 4457
$ws.Cells.Item(40, 6).Value = 276
$ws.Cells.Item(43, 6).Value = 2127
$ws.Cells.Item(47, 6).Value = 36
$ws.Cells.Item(48, 6).Value = 4305
$ws.Cells.Item(49, 6).Value = 197
